$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original style/format of column D so that setting text-like
# numeric strings (e.g. "582.90") does not get auto-converted to a number.
$colD = $ws.Range("D2:D51")
$origStyle = $colD.Style
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "62.305.76"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "3.007.25"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "582.90"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "145.68"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "3.003.61"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").Value = "34.47"
$ws.Range("E14").Value = "  -6.42%  "
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "3.501.00"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "7.11"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "62.278.15"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "3.008.60"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").Value = "456.57"
$ws.Range("E20").Value = "  -6.91%  "
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").Value = "7.40"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "81.79"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -10.33%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  -6.30%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").Value = "6.96"
$ws.Range("E31").Value = "  -6.24%  "
$ws.Range("D32").Value = "28.44"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("D35").Value = "0.0₃0798"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").Value = "5.76"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("D38").Value = "2.11"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("D39").Value = "9.17"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "50.24"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  -12.54%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").Value = "391.79"
$ws.Range("E43").Value = "  -10.36%  "
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  -7.68%  "
$ws.Range("D46").Value = "2.725.73"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "36.77"
$ws.Range("E47").Value = "  -6.52%  "
$ws.Range("D48").Value = "128.34"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("E51").Value = "  -2.03%  "

# Restore the original column style so no residual text-format override
# is left behind on the cells (matches the target workbook formatting).
$colD.Style = $origStyle

